$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, forcing text storage so that
# numeric-looking strings (e.g. "0.999", "549.87") are preserved verbatim
# instead of being parsed into floating point numbers by Excel.
function Set-TextValue {
    param($cell, [string]$value, $styleSourceCell)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $styleSourceCell.Style
}

Set-TextValue $ws.Range("D2") "63.374.07" $ws.Range("C2")
$ws.Range("E2").Value = "  -4.27%  "
Set-TextValue $ws.Range("D3") "3.095.69" $ws.Range("C3")
$ws.Range("E3").Value = "  -4.92%  "
Set-TextValue $ws.Range("D4") "0.999" $ws.Range("C4")
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue $ws.Range("D5") "549.87" $ws.Range("C5")
$ws.Range("E5").Value = "  -4.74%  "
Set-TextValue $ws.Range("D6") "137.42" $ws.Range("C6")
$ws.Range("E6").Value = "  -11.06%  "
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue $ws.Range("D8") "3.086.91" $ws.Range("C8")
$ws.Range("E8").Value = "  -4.97%  "
Set-TextValue $ws.Range("D9") "0.498" $ws.Range("C9")
$ws.Range("E9").Value = "  -3.66%  "
$ws.Range("E10").Value = "  -5.56%  "
$ws.Range("E11").Value = "  -11.22%  "
$ws.Range("E12").Value = "  -4.77%  "
Set-TextValue $ws.Range("D13") "35.61" $ws.Range("C13")
$ws.Range("E13").Value = "  -6.20%  "
$ws.Range("E14").Value = "  -8.35%  "
Set-TextValue $ws.Range("D15") "3.588.60" $ws.Range("C15")
$ws.Range("E15").Value = "  -4.94%  "
Set-TextValue $ws.Range("D16") "63.293.76" $ws.Range("C16")
$ws.Range("E16").Value = "  -4.54%  "
Set-TextValue $ws.Range("D17") "0.112" $ws.Range("C17")
$ws.Range("E17").Value = "  -3.21%  "
Set-TextValue $ws.Range("D18") "3.096.64" $ws.Range("C18")
$ws.Range("E18").Value = "  -4.80%  "
Set-TextValue $ws.Range("D19") "6.75" $ws.Range("C19")
$ws.Range("E19").Value = "  -5.87%  "
Set-TextValue $ws.Range("D20") "489.61" $ws.Range("C20")
$ws.Range("E20").Value = "  -12.78%  "
Set-TextValue $ws.Range("D21") "13.68" $ws.Range("C21")
$ws.Range("E21").Value = "  -5.89%  "
Set-TextValue $ws.Range("D22") "0.721" $ws.Range("C22")
$ws.Range("E22").Value = "  -3.60%  "
Set-TextValue $ws.Range("D23") "7.28" $ws.Range("C23")
$ws.Range("E23").Value = "  -7.44%  "
Set-TextValue $ws.Range("D24") "79.25" $ws.Range("C24")
$ws.Range("E24").Value = "  -3.99%  "
Set-TextValue $ws.Range("D25") "12.41" $ws.Range("C25")
$ws.Range("E25").Value = "  -9.22%  "
$ws.Range("E26").Value = "  -0.27%  "
Set-TextValue $ws.Range("D27") "8.48" $ws.Range("C27")
$ws.Range("E27").Value = "  -9.65%  "
Set-TextValue $ws.Range("D28") "2.76" $ws.Range("C28")
$ws.Range("E28").Value = "  -7.31%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D29") "1.98" $ws.Range("C29")
$ws.Range("E29").Value = "  -12.20%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D30") "1.00" $ws.Range("C30")
$ws.Range("E30").Value = "  -0.16%  "
Set-TextValue $ws.Range("D31") "26.65" $ws.Range("C31")
$ws.Range("E31").Value = "  -4.88%  "
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("E33").Value = "  -9.30%  "
Set-TextValue $ws.Range("D34") "57.93" $ws.Range("C34")
$ws.Range("E34").Value = "  +4.85%  "
Set-TextValue $ws.Range("D35") "515.29" $ws.Range("C35")
$ws.Range("E35").Value = "  -9.55%  "
Set-TextValue $ws.Range("D36") "6.03" $ws.Range("C36")
$ws.Range("E36").Value = "  -6.32%  "
Set-TextValue $ws.Range("D37") "5.15" $ws.Range("C37")
$ws.Range("E37").Value = "  -11.24%  "
$ws.Range("E38").Value = "  -12.87%  "
Set-TextValue $ws.Range("D39") "3.158.87" $ws.Range("C39")
$ws.Range("E39").Value = "  -0.21%  "
Set-TextValue $ws.Range("D40") "0.0805" $ws.Range("C40")
$ws.Range("E40").Value = "  -7.57%  "
Set-TextValue $ws.Range("D41") "0.120" $ws.Range("C41")
$ws.Range("E41").Value = "  -6.74%  "
Set-TextValue $ws.Range("D42") "8.18" $ws.Range("C42")
$ws.Range("E42").Value = "  -5.62%  "
Set-TextValue $ws.Range("D43") "2.67" $ws.Range("C43")
$ws.Range("E43").Value = "  -13.41%  "
Set-TextValue $ws.Range("D44") "0.260" $ws.Range("C44")
$ws.Range("E44").Value = "  -5.96%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D46") "2.07" $ws.Range("C46")
$ws.Range("E46").Value = "  -10.42%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "25.29" $ws.Range("C47")
$ws.Range("E47").Value = "  -5.47%  "
$ws.Range("E48").Value = "  -3.80%  "
$ws.Range("E49").Value = "  -4.39%  "
Set-TextValue $ws.Range("D50") "0.0₃0505" $ws.Range("C50")
$ws.Range("E50").Value = "  -9.66%  "
$ws.Range("E51").Value = "  -9.93%  "
